# "feat: Implemented System Login"
#
# The Settings sheet's System1 section is reshuffled: the blank spacer rows
# (9, 11, 13) are removed so the three existing entries (SHA1_URL,
# ExceptionEmail, System1_Credential) move up two rows each, and a brand new
# "OrchestratorFolder" row (with an explanatory description in column C) is
# added where the old trailing blank rows (16-18) used to be. The two
# external hyperlinks move together with the cells that now hold their
# display text, and the description column is widened to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")
$ws.Activate()

# --- Row 10: SHA1_URL / sha1.html link -> moves to row 12, loses its old content ---
# --- Row 12: ExceptionEmail / mailto link -> moves to row 14 ---
# --- Row 14: System1Credential / ACME_Credential -> becomes System1_Credential at row 10 ---
# Work from the bottom up so we never overwrite data we still need to read.

$ws.Range("A14").Value = "ExceptionEmail"
$ws.Range("B14").Value = "exceptions@acme-test.com"

$ws.Range("A12").Value = "SHA1_URL"
$ws.Range("B12").Value = "https://emn178.github.io/online-tools/sha1.html"

$ws.Range("A10").Value = "System1_Credential"
$ws.Range("B10").Value = "ACME_Credential"
$ws.Range("B10").Style = "Normal"

# --- New row 16: OrchestratorFolder / ACME_Automation / description ---
$ws.Range("A16").Value = "OrchestratorFolder"
$ws.Range("B16").Value = "ACME_Automation"
$ws.Range("C16").Value = "Folder name. The value must match a folder defined in Orchestrator. For classic folders leave the value field empty."

# --- Rebuild the two external hyperlinks against their new cells ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B14"), "mailto:exceptions@acme-test.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://emn178.github.io/online-tools/sha1.html") | Out-Null

# Hyperlinks.Add re-stamps its own cell format; re-apply the named
# "Hyperlink" style explicitly afterwards so both cells land back on the
# workbook's existing Hyperlink cell style instead of a freshly minted one.
$ws.Range("B14").Style = "Hyperlink"
$ws.Range("B12").Style = "Hyperlink"

# --- Widen column C to fit the new, much longer description text ---
$ws.Columns.Item(3).ColumnWidth = 179.3333

# --- Update the saved selection to match the author's final cursor position ---
$ws.Range("B24").Select() | Out-Null
